$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.931.27"
$ws.Range("E2").Value = "  -3.75%  "

$ws.Range("D3").Value = "1.637.86"
$ws.Range("E3").Value = "  -5.97%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9972"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.35%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.30%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4712"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -6.84%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2559"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.83%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06010"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.69%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07135"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.35%  "

$ws.Range("D11").Value = "1.633.65"
$ws.Range("E11").Value = "  -6.19%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.83"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.45%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6150"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.91%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.406"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.72%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "72.65"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.20%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.09%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9975"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.33%  "

$ws.Range("D18").Value = "24.919.10"
$ws.Range("E18").Value = "  -3.83%  "

$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.88%  "

$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006563"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.53%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.399"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.82%  "

$ws.Range("D22").Value = "1.844.30"
$ws.Range("E22").Value = "  -6.23%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.563"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.68%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.257"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.09%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "132.78"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.53%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.83"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.47%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.372"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -8.77%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "102.59"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.60%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.655"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.14%  "

$ws.Range("E30").Value = "  -4.73%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07733"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.99%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.546"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.17%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04364"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.52%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9986"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.13%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.596"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.19%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9204"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.26%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5803"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.34%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.539"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.77%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01554"
$ws.Range("D39").Style = "Normal"

$ws.Range("E40").Value = "  -0.26%  "

$ws.Range("B41").Value = "PaxosStandard"
$ws.Range("C41").Value = "https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9985"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.32%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8130"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.26%  "

$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.800"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.66%  "

$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "97.50"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.44%  "

$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3702"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.91%  "

$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.729"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.99%  "

$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1124"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.19%  "

$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05225"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.15%  "

$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.077"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.80%  "

$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "29.49"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.52%  "

$ws.Range("B51").Value = "TrueUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9992"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.37%  "
